$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"7.985720451708375e-18"
$ws.Range("C2").Value = [double]"0.0001000769650044792"
$ws.Range("D2").Value = [double]"0.6683886250799669"
$ws.Range("E2").Value = [double]"0.9047753617887724"
$ws.Range("F2").Value = [double]"0.8327056364760048"
$ws.Range("G2").Value = [double]"0.5565824650942856"
$ws.Range("H2").Value = [double]"0.7676176443735369"
$ws.Range("I2").Value = [double]"0.992326740186214"
$ws.Range("J2").Value = [double]"0.2427711145223551"
$ws.Range("K2").Value = [double]"0.4060235651730132"
$ws.Range("L2").Value = [double]"0.4447740396687111"
$ws.Range("M2").Value = [double]"0.4793335286592828"
$ws.Range("N2").Value = [double]"0.5187882833856901"
$ws.Range("O2").Value = [double]"0.4810630054064726"
$ws.Range("P2").Value = [double]"0.330429889411633"
$ws.Range("Q2").Value = [double]"0.4316527563557564"
$ws.Range("R2").Value = [double]"0.5058178246372492"
$ws.Range("S2").Value = [double]"0.4681871046562327"
$ws.Range("T2").Value = [double]"0.04729300264716872"
$ws.Range("U2").Value = [double]"8.887395411605914e-19"
$ws.Range("V2").Value = [double]"4.139161972636311e-11"
$ws.Range("B3").Value = [double]"2.881059432383614e-20"
$ws.Range("C3").Value = [double]"9.192938627272448e-07"
$ws.Range("D3").Value = [double]"0.9044543795128401"
$ws.Range("E3").Value = [double]"0.6857116125642075"
$ws.Range("F3").Value = [double]"0.5032834514357978"
$ws.Range("G3").Value = [double]"0.800260162688039"
$ws.Range("H3").Value = [double]"0.9688444758471388"
$ws.Range("I3").Value = [double]"0.7652990112723463"
$ws.Range("J3").Value = [double]"0.1133509186535721"
$ws.Range("K3").Value = [double]"0.1579812195996347"
$ws.Range("L3").Value = [double]"0.2225295456529744"
$ws.Range("M3").Value = [double]"0.1738957980271216"
$ws.Range("N3").Value = [double]"0.2539441463819709"
$ws.Range("O3").Value = [double]"0.1929912518661896"
$ws.Range("P3").Value = [double]"0.165405359663531"
$ws.Range("Q3").Value = [double]"0.2353433340135165"
$ws.Range("R3").Value = [double]"0.1906973973297687"
$ws.Range("S3").Value = [double]"0.2149738935988811"
$ws.Range("T3").Value = [double]"0.005206012867619276"
$ws.Range("U3").Value = [double]"5.040322510300809e-33"
$ws.Range("V3").Value = [double]"2.004935471430351e-13"
$ws.Range("B4").Value = [double]"5.537709910156304e-25"
$ws.Range("C4").Value = [double]"2.524652071849339e-11"
$ws.Range("D4").Value = [double]"0.8770672689597531"
$ws.Range("E4").Value = [double]"0.4821106777989222"
$ws.Range("F4").Value = [double]"0.3917216472188526"
$ws.Range("G4").Value = [double]"0.9245286878763932"
$ws.Range("H4").Value = [double]"0.8387633316278917"
$ws.Range("I4").Value = [double]"0.6712692810110672"
$ws.Range("J4").Value = [double]"0.08539621490894032"
$ws.Range("K4").Value = [double]"0.08589747403552649"
$ws.Range("L4").Value = [double]"0.1354605806541475"
$ws.Range("M4").Value = [double]"0.1042754438970071"
$ws.Range("N4").Value = [double]"0.1660910053643732"
$ws.Range("O4").Value = [double]"0.1072282228019529"
$ws.Range("P4").Value = [double]"0.1308820165495146"
$ws.Range("Q4").Value = [double]"0.2067485153718752"
$ws.Range("R4").Value = [double]"0.1341457576953057"
$ws.Range("S4").Value = [double]"0.1113037170933485"
$ws.Range("T4").Value = [double]"0.001783251267607191"
$ws.Range("U4").Value = [double]"8.270574016290035e-49"
$ws.Range("V4").Value = [double]"8.903939292390426e-15"
$ws.Range("B5").Value = [double]"0.09093521525652967"
$ws.Range("C5").Value = [double]"0.001162627515801274"
$ws.Range("D5").Value = [double]"0.9645388028330022"
$ws.Range("E5").Value = [double]"0.6322876399121994"
$ws.Range("F5").Value = [double]"0.4113855118694248"
$ws.Range("G5").Value = [double]"0.9669684514885055"
$ws.Range("H5").Value = [double]"0.6993044446589063"
$ws.Range("I5").Value = [double]"0.7607992781442277"
$ws.Range("J5").Value = [double]"0.7844244945315254"
$ws.Range("K5").Value = [double]"0.531078556748204"
$ws.Range("L5").Value = [double]"0.7687906264273019"
$ws.Range("M5").Value = [double]"0.9088851595551664"
$ws.Range("N5").Value = [double]"0.9149070619196777"
$ws.Range("O5").Value = [double]"0.3426666060866215"
$ws.Range("P5").Value = [double]"0.9038562479833077"
$ws.Range("Q5").Value = [double]"0.8635179134628745"
$ws.Range("R5").Value = [double]"0.895576053185705"
$ws.Range("S5").Value = [double]"0.693035700587866"
$ws.Range("T5").Value = [double]"0.4204763185495483"
$ws.Range("U5").Value = [double]"1.300700635717633e-20"
$ws.Range("V5").Value = [double]"8.669010199122036e-53"
